# Training update to confirmed barriers tracking table
#
# Fills in a worked/training example in row 2 of the "Confirmed barriers"
# sheet (one value per tracked column) and leaves the selection on the last
# column of that row, mirroring the edit made in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmed barriers")

$ws.Range("A2").Value = "Test 1"                              # Internal Name
$ws.Range("B2").Value = 123                                   # Barrier ID (Crossing ID)
$ws.Range("C2").Value = "Stream"                               # Watercourse name
$ws.Range("D2").Value = "Road"                                 # Road name
$ws.Range("H2").Value = "Dam"                                   # Barrier type
$ws.Range("I2").Value = "Mine"                                  # Barrier Owner
$ws.Range("J2").Value = "Detailed habitat investigation"        # Assessment Type/Assessment Step Completed
$ws.Range("L2").Value = "Medium"                                # Upstream habitat Quality
$ws.Range("M2").Value = "Moderate"                              # Constructability
$ws.Range("N2").Value = 5000                                    # Estimated cost
$ws.Range("P2").Value = "High"                                  # Priority
$ws.Range("Q2").Value = "Commission engineering designs"        # Next Steps
$ws.Range("R2").Value = 2025                                    # Timeline for Next Steps
$ws.Range("S2").Value = "CWF"                                   # Lead for Next Steps
$ws.Range("T2").Value = "LDN"                                   # Others Involved in Next Steps
$ws.Range("U2").Value = "Test."                                 # Reason
$ws.Range("V2").Value = "Test."                                 # Notes

# Match the cursor position left after the edit.
$ws.Activate()
$ws.Range("V2").Select()
